$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1348.1384
$ws.Range("I15").Value = 1348.1384
$ws.Range("K15").Value = 4044.4152
$ws.Range("M15").Value = -3875.4152
$ws.Range("H42").Value = 1389289.6
$ws.Range("I42").Value = 5000211
$ws.Range("J42").Value = 473.69232
$ws.Range("K42").Value = 15000633
$ws.Range("L42").Value = 1421.07696
$ws.Range("M42").Value = -15000403
$ws.Range("N42").Value = -1881.07696
$ws.Range("H82").Value = 2899.9333
$ws.Range("I82").Value = 699.8
$ws.Range("K82").Value = 2099.4
$ws.Range("M82").Value = -1693.4
$ws.Range("H85").Value = 2899.9333
$ws.Range("I85").Value = 699.8
$ws.Range("K85").Value = 2099.4
$ws.Range("M85").Value = -695.3999999999996
$ws.Range("H113").Value = 168667.5
$ws.Range("I113").Value = 202001
$ws.Range("K113").Value = 202001
$ws.Range("M113").Value = -198747
$ws.Range("H129").Value = 942.0213
$ws.Range("J129").Value = 985.0732
$ws.Range("L129").Value = 2955.2196
$ws.Range("N129").Value = -12955.2196
$ws.Range("H132").Value = 5957705.5
$ws.Range("I132").Value = 6762667
$ws.Range("J132").Value = 988
$ws.Range("K132").Value = 20288001
$ws.Range("L132").Value = 2964
$ws.Range("M132").Value = -20285471
$ws.Range("N132").Value = -8024
$ws.Range("H138").Value = 4296.085
$ws.Range("I138").Value = 2151.2727
$ws.Range("J138").Value = 5571.3784
$ws.Range("K138").Value = 6453.8181
$ws.Range("L138").Value = 16714.1352
$ws.Range("M138").Value = -1313.8181
$ws.Range("N138").Value = -26994.1352
$ws.Range("H141").Value = 2685.4783
$ws.Range("I141").Value = 1798.1428
$ws.Range("J141").Value = 12002.5
$ws.Range("K141").Value = 5394.428400000001
$ws.Range("L141").Value = 36007.5
$ws.Range("M141").Value = -214.4284000000007
$ws.Range("N141").Value = -46367.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23390.025
$ws.Range("I32").Value = 3855.9697
$ws.Range("J32").Value = 130827.336
$ws.Range("K32").Value = 3855.9697
$ws.Range("L32").Value = 130827.336
$ws.Range("M32").Value = -3568.9697
$ws.Range("N32").Value = -131401.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16770.533
$ws.Range("I82").Value = 3936.889
$ws.Range("J82").Value = 36021
$ws.Range("K82").Value = 3936.889
$ws.Range("L82").Value = 36021
$ws.Range("M82").Value = -3553.889
$ws.Range("N82").Value = -36787
$ws.Range("H85").Value = 16770.533
$ws.Range("I85").Value = 3936.889
$ws.Range("J85").Value = 36021
$ws.Range("K85").Value = 3936.889
$ws.Range("L85").Value = 36021
$ws.Range("M85").Value = -2610.889
$ws.Range("N85").Value = -38673
$ws.Range("H86").Value = 50006.914
$ws.Range("I86").Value = 75204.87
$ws.Range("J86").Value = 2760.75
$ws.Range("K86").Value = 75204.87
$ws.Range("L86").Value = 2760.75
$ws.Range("M86").Value = -74081.87
$ws.Range("N86").Value = -5006.75
$ws.Range("H89").Value = 50006.914
$ws.Range("I89").Value = 75204.87
$ws.Range("J89").Value = 2760.75
$ws.Range("K89").Value = 376024.35
$ws.Range("L89").Value = 13803.75
$ws.Range("M89").Value = -370408.35
$ws.Range("N89").Value = -25035.75
$ws.Range("H105").Value = 78859.195
$ws.Range("I105").Value = 51816.4
$ws.Range("J105").Value = 169001.83
$ws.Range("K105").Value = 51816.4
$ws.Range("L105").Value = 169001.83
$ws.Range("M105").Value = -50069.4
$ws.Range("N105").Value = -172495.83
$ws.Range("H134").Value = 7802.4
$ws.Range("I134").Value = 8753
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 26259
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -23724
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 48250
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H24").Value = 48250
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H31").Value = 15076.333
$ws.Range("I31").Value = 32243.844
$ws.Range("J31").Value = 2300.5117
$ws.Range("K31").Value = 32243.844
$ws.Range("L31").Value = 2300.5117
$ws.Range("M31").Value = -31948.844
$ws.Range("N31").Value = -2890.5117
$ws.Range("H34").Value = 15076.333
$ws.Range("I34").Value = 32243.844
$ws.Range("J34").Value = 2300.5117
$ws.Range("K34").Value = 32243.844
$ws.Range("L34").Value = 2300.5117
$ws.Range("M34").Value = -32041.844
$ws.Range("N34").Value = -2704.5117
$ws.Range("H58").Value = 10385.241
$ws.Range("I58").Value = 1691.55
$ws.Range("J58").Value = 29704.555
$ws.Range("K58").Value = 1691.55
$ws.Range("L58").Value = 29704.555
$ws.Range("M58").Value = -1488.55
$ws.Range("N58").Value = -30110.555
$ws.Range("H99").Value = 9892.875
$ws.Range("I99").Value = 3986.6667
$ws.Range("J99").Value = 17486.572
$ws.Range("K99").Value = 3986.6667
$ws.Range("L99").Value = 17486.572
$ws.Range("M99").Value = -2488.6667
$ws.Range("N99").Value = -20482.572
$ws.Range("H107").Value = 4379.185
$ws.Range("I107").Value = 7828.9287
$ws.Range("J107").Value = 664.0769
$ws.Range("K107").Value = 7828.9287
$ws.Range("L107").Value = 664.0769
$ws.Range("M107").Value = -5908.9287
$ws.Range("N107").Value = -4504.0769
$ws.Range("H126").Value = 9892.875
$ws.Range("I126").Value = 3986.6667
$ws.Range("J126").Value = 17486.572
$ws.Range("K126").Value = 11960.0001
$ws.Range("L126").Value = 52459.716
$ws.Range("M126").Value = -9490.000100000001
$ws.Range("N126").Value = -57399.716
$ws.Range("H132").Value = 4288.2
$ws.Range("I132").Value = 4611.1
$ws.Range("J132").Value = 3642.4
$ws.Range("K132").Value = 13833.3
$ws.Range("L132").Value = 10927.2
$ws.Range("M132").Value = -11303.3
$ws.Range("N132").Value = -15987.2
$ws.Range("H134").Value = 1631.5
$ws.Range("I134").Value = 1292.8823
$ws.Range("J134").Value = 2782.8
$ws.Range("K134").Value = 3878.6469
$ws.Range("L134").Value = 8348.400000000001
$ws.Range("M134").Value = -1343.6469
$ws.Range("N134").Value = -13418.4
$ws.Range("H136").Value = 10385.241
$ws.Range("I136").Value = 1691.55
$ws.Range("J136").Value = 29704.555
$ws.Range("K136").Value = 5074.65
$ws.Range("L136").Value = 89113.66500000001
$ws.Range("M136").Value = -2524.65
$ws.Range("N136").Value = -94213.66500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 505.64285
$ws.Range("I26").Value = 166.66667
$ws.Range("J26").Value = 759.875
$ws.Range("K26").Value = 500.00001
$ws.Range("L26").Value = 2279.625
$ws.Range("M26").Value = -212.00001
$ws.Range("N26").Value = -2855.625
$ws.Range("H122").Value = 6353.3887
$ws.Range("I122").Value = 272.8
$ws.Range("J122").Value = 8692.077
$ws.Range("K122").Value = 2455.2
$ws.Range("L122").Value = 78228.693
$ws.Range("M122").Value = -5.200000000000273
$ws.Range("N122").Value = -83128.693
$ws.Range("H131").Value = 1548.7954
$ws.Range("I131").Value = 1872.5
$ws.Range("J131").Value = 1516.425
$ws.Range("K131").Value = 5617.5
$ws.Range("L131").Value = 4549.275
$ws.Range("M131").Value = -577.5
$ws.Range("N131").Value = -14629.275
$ws.Range("H137").Value = 17546394
$ws.Range("J137").Value = 37040064
$ws.Range("L137").Value = 111120192
$ws.Range("N137").Value = -111130392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 76926050
$ws.Range("I97").Value = 125003490
$ws.Range("J97").Value = 2142
$ws.Range("K97").Value = 125003490
$ws.Range("L97").Value = 2142
$ws.Range("M97").Value = -125002994
$ws.Range("N97").Value = -3134
$ws.Range("H132").Value = 5033.6665
$ws.Range("I132").Value = 5023.5
$ws.Range("J132").Value = 5054
$ws.Range("K132").Value = 15070.5
$ws.Range("L132").Value = 15162
$ws.Range("M132").Value = -12540.5
$ws.Range("N132").Value = -20222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 7501
$ws.Range("I21").Value = 1006
$ws.Range("J21").Value = 8800
$ws.Range("K21").Value = 1006
$ws.Range("L21").Value = 8800
$ws.Range("M21").Value = -832
$ws.Range("N21").Value = -9148
$ws.Range("H40").Value = 168630.67
$ws.Range("I40").Value = 251076
$ws.Range("J40").Value = 3740
$ws.Range("K40").Value = 251076
$ws.Range("L40").Value = 3740
$ws.Range("M40").Value = -250940
$ws.Range("N40").Value = -4012
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 500000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H106").Value = 18345
$ws.Range("J106").Value = 18345
$ws.Range("L106").Value = 18345
$ws.Range("N106").Value = -20869
$ws.Range("H126").Value = 3756.8
$ws.Range("I126").Value = 3946
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 11838
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -9368
$ws.Range("N126").Value = -13940
$ws.Range("H127").Value = 29995
$ws.Range("J127").Value = 29995
$ws.Range("L127").Value = 29995
$ws.Range("N127").Value = -39915
$ws.Range("H132").Value = 12272
$ws.Range("I132").Value = 15908.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 47725.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -45195.5
$ws.Range("N132").Value = -20057
